$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.910.53'
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = '1.641.48'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -1.56%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.13'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5053'
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("E7").Value = '  -0.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2579'
$ws.Range("E8").Value = '  +0.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06406'
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.56'
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07783'
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").Value = '1.651.30'
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.268'
$ws.Range("E13").Value = '  +0.86%  '
$ws.Range("D14").Value = '1.867.87'
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5444'
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").Value = '0.0₅7927'
$ws.Range("E16").Value = '  +0.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.62'
$ws.Range("E17").Value = '  +2.51%  '
$ws.Range("D18").Value = '25.958.47'
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  -1.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '198.36'
$ws.Range("E20").Value = '  -2.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.384'
$ws.Range("E21").Value = '  +1.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.907'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.974'
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.006'
$ws.Range("E24").Value = '  -0.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.877'
$ws.Range("E25").Value = '  -4.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.52'
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1136'
$ws.Range("E27").Value = '  -1.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.831'
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.72'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.241'
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04945'
$ws.Range("E31").Value = '  -0.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.265'
$ws.Range("E32").Value = '  +0.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.202'
$ws.Range("E33").Value = '  +0.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.540'
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.372'
$ws.Range("E35").Value = '  +1.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8939'
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.608'
$ws.Range("E37").Value = '  -1.35%  '
$ws.Range("D38").Value = '1.144.54'
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5560'
$ws.Range("E39").Value = '  -1.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01570'
$ws.Range("E40").Value = '  +1.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.006'
$ws.Range("E41").Value = '  -0.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.729'
$ws.Range("E42").Value = '  +1.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8123'
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.82'
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("E45").Value = '  +7.07%  '
$ws.Range("D46").Value = '1.778.83'
$ws.Range("E46").Value = '  +0.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4511'
$ws.Range("E47").Value = '  -0.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.005'
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.80'
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05059'
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.005'
